# Add a new "2021" column (R) that mirrors the existing "2020" column (Q):
# same per-row formatting, plus the new data values/formulas for that year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone column Q's formatting onto the new column R for the data block (rows 4-44).
$ws.Range("Q4:Q44").Copy()
$ws.Range("R4:R44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Populate the 2021 header and data values/formulas in column R.
$ws.Range("R4").Value = 2021

$ws.Range("R6").Formula = "=R7/R8*100"
$ws.Range("R7").Value = 1931.83
$ws.Range("R8").Value = 739818.5

$ws.Range("R10").Formula = "=R11/R12*100"
$ws.Range("R11").Value = 1552.9
$ws.Range("R12").Value = 25048.6

$ws.Range("R14").Formula = "=R15/R16*100"
$ws.Range("R15").Value = 125.7
$ws.Range("R16").Value = 82213.9

$ws.Range("R18").Formula = "=R19/R20*100"
$ws.Range("R19").Value = 99.6
$ws.Range("R20").Value = 80059.6

$ws.Range("R22").Formula = "=R23/R24*100"
$ws.Range("R23").Value = 0.9
$ws.Range("R24").Value = 17172.7

$ws.Range("R26").Formula = "=R27/R28*100"
$ws.Range("R27").Value = 15.9
$ws.Range("R28").Value = 56666.5

$ws.Range("R30").Formula = "=R31/R32*100"
$ws.Range("R31").Value = 58.5
$ws.Range("R32").Value = 30765.1

$ws.Range("R34").Formula = "=R35/R36*100"
$ws.Range("R35").Value = 78.3
$ws.Range("R36").Value = 110267.1

$ws.Range("R38").Value = "-"
$ws.Range("R39").Value = "-"
$ws.Range("R40").Value = 297797.6

$ws.Range("R42").Value = "-"
$ws.Range("R43").Value = "-"
$ws.Range("R44").Value = 39827.4

# 3) Update the active selection to match where the author ended up after adding the column.
$ws.Range("T9").Select()
